$wb = $excel.ActiveWorkbook

# --- Sheet1: LP1912 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A2").Value = "Última actualización: 16:45:34"
$ws1.Range("A3").Value = "Total filas: 336"

# Reassign shuffled tie-break rows (values swapped due to upstream re-sort)
$ws1.Cells.Item(120,1).Value = "10:13:53"
$ws1.Cells.Item(120,2).Value = "10:34"
$ws1.Cells.Item(120,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(120,4).Value = 21
$ws1.Cells.Item(120,5).Value = "LP1912"
$ws1.Cells.Item(121,1).Value = "10:13:53"
$ws1.Cells.Item(121,2).Value = "10:34"
$ws1.Cells.Item(121,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(121,4).Value = 21
$ws1.Cells.Item(121,5).Value = "LP1912"
$ws1.Cells.Item(162,1).Value = "10:13:53"
$ws1.Cells.Item(162,2).Value = "12:06"
$ws1.Cells.Item(162,3).Value = "14_ABASTO"
$ws1.Cells.Item(162,4).Value = 113
$ws1.Cells.Item(162,5).Value = "LP1912"
$ws1.Cells.Item(164,1).Value = "10:13:53"
$ws1.Cells.Item(164,2).Value = "12:06"
$ws1.Cells.Item(164,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(164,4).Value = 113
$ws1.Cells.Item(164,5).Value = "LP1912"
$ws1.Cells.Item(173,1).Value = "12:01:11"
$ws1.Cells.Item(173,2).Value = "12:21"
$ws1.Cells.Item(173,3).Value = "215A_EL PATO"
$ws1.Cells.Item(173,4).Value = 20
$ws1.Cells.Item(173,5).Value = "LP1912"
$ws1.Cells.Item(174,1).Value = "12:01:11"
$ws1.Cells.Item(174,2).Value = "12:21"
$ws1.Cells.Item(174,3).Value = "14_ABASTO"
$ws1.Cells.Item(174,4).Value = 20
$ws1.Cells.Item(174,5).Value = "LP1912"
$ws1.Cells.Item(175,1).Value = "10:52:37"
$ws1.Cells.Item(175,2).Value = "12:21"
$ws1.Cells.Item(175,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(175,4).Value = 89
$ws1.Cells.Item(175,5).Value = "LP1912"
$ws1.Cells.Item(184,1).Value = "12:35:30"
$ws1.Cells.Item(184,2).Value = "12:38"
$ws1.Cells.Item(184,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(184,4).Value = 3
$ws1.Cells.Item(184,5).Value = "LP1912"
$ws1.Cells.Item(185,1).Value = "10:52:37"
$ws1.Cells.Item(185,2).Value = "12:38"
$ws1.Cells.Item(185,3).Value = "17_179 Y 38"
$ws1.Cells.Item(185,4).Value = 106
$ws1.Cells.Item(185,5).Value = "LP1912"
$ws1.Cells.Item(193,1).Value = "12:50:41"
$ws1.Cells.Item(193,2).Value = "12:50"
$ws1.Cells.Item(193,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(193,4).Value = 0
$ws1.Cells.Item(193,5).Value = "LP1912"
$ws1.Cells.Item(194,1).Value = "10:52:37"
$ws1.Cells.Item(194,2).Value = "12:50"
$ws1.Cells.Item(194,3).Value = "15_ABASTO"
$ws1.Cells.Item(194,4).Value = 118
$ws1.Cells.Item(194,5).Value = "LP1912"
$ws1.Cells.Item(195,1).Value = "12:01:11"
$ws1.Cells.Item(195,2).Value = "13:02"
$ws1.Cells.Item(195,3).Value = "15_ABASTO"
$ws1.Cells.Item(195,4).Value = 61
$ws1.Cells.Item(195,5).Value = "LP1912"
$ws1.Cells.Item(196,1).Value = "12:35:30"
$ws1.Cells.Item(196,2).Value = "13:02"
$ws1.Cells.Item(196,3).Value = "14_ABASTO"
$ws1.Cells.Item(196,4).Value = 27
$ws1.Cells.Item(196,5).Value = "LP1912"
$ws1.Cells.Item(213,1).Value = "13:18:40"
$ws1.Cells.Item(213,2).Value = "13:36"
$ws1.Cells.Item(213,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(213,4).Value = 18
$ws1.Cells.Item(213,5).Value = "LP1912"
$ws1.Cells.Item(214,1).Value = "13:18:40"
$ws1.Cells.Item(214,2).Value = "13:36"
$ws1.Cells.Item(214,3).Value = "15_ABASTO"
$ws1.Cells.Item(214,4).Value = 18
$ws1.Cells.Item(214,5).Value = "LP1912"
$ws1.Cells.Item(216,1).Value = "13:18:40"
$ws1.Cells.Item(216,2).Value = "13:46"
$ws1.Cells.Item(216,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(216,4).Value = 28
$ws1.Cells.Item(216,5).Value = "LP1912"
$ws1.Cells.Item(217,1).Value = "12:01:11"
$ws1.Cells.Item(217,2).Value = "13:46"
$ws1.Cells.Item(217,3).Value = "17_ROMERO"
$ws1.Cells.Item(217,4).Value = 105
$ws1.Cells.Item(217,5).Value = "LP1912"
$ws1.Cells.Item(305,1).Value = "16:45:34"
$ws1.Cells.Item(305,2).Value = "16:47"
$ws1.Cells.Item(305,3).Value = "14_ABASTO"
$ws1.Cells.Item(305,4).Value = 2
$ws1.Cells.Item(305,5).Value = "LP1912"
$ws1.Cells.Item(306,1).Value = "15:36:13"
$ws1.Cells.Item(306,2).Value = "16:48"
$ws1.Cells.Item(306,3).Value = "15_ABASTO"
$ws1.Cells.Item(306,4).Value = 72
$ws1.Cells.Item(306,5).Value = "LP1912"
$ws1.Cells.Item(307,1).Value = "16:34:19"
$ws1.Cells.Item(307,2).Value = "16:50"
$ws1.Cells.Item(307,3).Value = "14_ABASTO"
$ws1.Cells.Item(307,4).Value = 16
$ws1.Cells.Item(307,5).Value = "LP1912"
$ws1.Cells.Item(308,1).Value = "15:59:02"
$ws1.Cells.Item(308,2).Value = "16:51"
$ws1.Cells.Item(308,3).Value = "14_ABASTO"
$ws1.Cells.Item(308,4).Value = 52
$ws1.Cells.Item(308,5).Value = "LP1912"
$ws1.Cells.Item(309,1).Value = "14:59:23"
$ws1.Cells.Item(309,2).Value = "16:56"
$ws1.Cells.Item(309,3).Value = "17_179 Y 38"
$ws1.Cells.Item(309,4).Value = 117
$ws1.Cells.Item(309,5).Value = "LP1912"
$ws1.Cells.Item(310,1).Value = "15:59:02"
$ws1.Cells.Item(310,2).Value = "16:57"
$ws1.Cells.Item(310,3).Value = "10_OLMOS"
$ws1.Cells.Item(310,4).Value = 58
$ws1.Cells.Item(310,5).Value = "LP1912"
$ws1.Cells.Item(311,1).Value = "16:34:19"
$ws1.Cells.Item(311,2).Value = "17:04"
$ws1.Cells.Item(311,3).Value = "215A_EL PATO"
$ws1.Cells.Item(311,4).Value = 30
$ws1.Cells.Item(311,5).Value = "LP1912"
$ws1.Cells.Item(312,1).Value = "16:45:34"
$ws1.Cells.Item(312,2).Value = "17:04"
$ws1.Cells.Item(312,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(312,4).Value = 19
$ws1.Cells.Item(312,5).Value = "LP1912"
$ws1.Cells.Item(313,1).Value = "16:45:34"
$ws1.Cells.Item(313,2).Value = "17:04"
$ws1.Cells.Item(313,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(313,4).Value = 19
$ws1.Cells.Item(313,5).Value = "LP1912"
$ws1.Cells.Item(314,1).Value = "16:20:15"
$ws1.Cells.Item(314,2).Value = "17:05"
$ws1.Cells.Item(314,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(314,4).Value = 45
$ws1.Cells.Item(314,5).Value = "LP1912"
$ws1.Cells.Item(315,1).Value = "15:36:13"
$ws1.Cells.Item(315,2).Value = "17:05"
$ws1.Cells.Item(315,3).Value = "215A_EL PATO"
$ws1.Cells.Item(315,4).Value = 89
$ws1.Cells.Item(315,5).Value = "LP1912"
$ws1.Cells.Item(316,1).Value = "16:34:19"
$ws1.Cells.Item(316,2).Value = "17:10"
$ws1.Cells.Item(316,3).Value = "10_OLMOS"
$ws1.Cells.Item(316,4).Value = 36
$ws1.Cells.Item(316,5).Value = "LP1912"
$ws1.Cells.Item(317,1).Value = "16:34:19"
$ws1.Cells.Item(317,2).Value = "17:16"
$ws1.Cells.Item(317,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(317,4).Value = 42
$ws1.Cells.Item(317,5).Value = "LP1912"
$ws1.Cells.Item(318,1).Value = "15:59:02"
$ws1.Cells.Item(318,2).Value = "17:17"
$ws1.Cells.Item(318,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(318,4).Value = 78
$ws1.Cells.Item(318,5).Value = "LP1912"
$ws1.Cells.Item(319,1).Value = "16:45:34"
$ws1.Cells.Item(319,2).Value = "17:20"
$ws1.Cells.Item(319,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(319,4).Value = 35
$ws1.Cells.Item(319,5).Value = "LP1912"
$ws1.Cells.Item(320,1).Value = "15:36:13"
$ws1.Cells.Item(320,2).Value = "17:21"
$ws1.Cells.Item(320,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(320,4).Value = 105
$ws1.Cells.Item(320,5).Value = "LP1912"
$ws1.Cells.Item(321,1).Value = "16:20:15"
$ws1.Cells.Item(321,2).Value = "17:21"
$ws1.Cells.Item(321,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(321,4).Value = 61
$ws1.Cells.Item(321,5).Value = "LP1912"
$ws1.Cells.Item(322,1).Value = "15:36:13"
$ws1.Cells.Item(322,2).Value = "17:24"
$ws1.Cells.Item(322,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(322,4).Value = 108
$ws1.Cells.Item(322,5).Value = "LP1912"
$ws1.Cells.Item(323,1).Value = "16:34:19"
$ws1.Cells.Item(323,2).Value = "17:28"
$ws1.Cells.Item(323,3).Value = "14_ABASTO"
$ws1.Cells.Item(323,4).Value = 54
$ws1.Cells.Item(323,5).Value = "LP1912"
$ws1.Cells.Item(324,1).Value = "16:34:19"
$ws1.Cells.Item(324,2).Value = "17:31"
$ws1.Cells.Item(324,3).Value = "15_ABASTO"
$ws1.Cells.Item(324,4).Value = 57
$ws1.Cells.Item(324,5).Value = "LP1912"
$ws1.Cells.Item(325,1).Value = "16:45:34"
$ws1.Cells.Item(325,2).Value = "17:34"
$ws1.Cells.Item(325,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(325,4).Value = 49
$ws1.Cells.Item(325,5).Value = "LP1912"
$ws1.Cells.Item(326,1).Value = "16:20:15"
$ws1.Cells.Item(326,2).Value = "17:36"
$ws1.Cells.Item(326,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(326,4).Value = 76
$ws1.Cells.Item(326,5).Value = "LP1912"
$ws1.Cells.Item(327,1).Value = "15:59:02"
$ws1.Cells.Item(327,2).Value = "17:37"
$ws1.Cells.Item(327,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(327,4).Value = 98
$ws1.Cells.Item(327,5).Value = "LP1912"
$ws1.Cells.Item(328,1).Value = "15:59:02"
$ws1.Cells.Item(328,2).Value = "17:38"
$ws1.Cells.Item(328,3).Value = "17_ROMERO"
$ws1.Cells.Item(328,4).Value = 99
$ws1.Cells.Item(328,5).Value = "LP1912"
$ws1.Cells.Item(329,1).Value = "16:45:34"
$ws1.Cells.Item(329,2).Value = "17:38"
$ws1.Cells.Item(329,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(329,4).Value = 53
$ws1.Cells.Item(329,5).Value = "LP1912"
$ws1.Cells.Item(330,1).Value = "16:34:19"
$ws1.Cells.Item(330,2).Value = "17:39"
$ws1.Cells.Item(330,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(330,4).Value = 65
$ws1.Cells.Item(330,5).Value = "LP1912"
$ws1.Cells.Item(331,1).Value = "16:45:34"
$ws1.Cells.Item(331,2).Value = "17:40"
$ws1.Cells.Item(331,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(331,4).Value = 55
$ws1.Cells.Item(331,5).Value = "LP1912"
$ws1.Cells.Item(332,1).Value = "15:59:02"
$ws1.Cells.Item(332,2).Value = "17:40"
$ws1.Cells.Item(332,3).Value = "215B_EL PATO"
$ws1.Cells.Item(332,4).Value = 101
$ws1.Cells.Item(332,5).Value = "LP1912"
$ws1.Cells.Item(333,1).Value = "16:34:19"
$ws1.Cells.Item(333,2).Value = "17:41"
$ws1.Cells.Item(333,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(333,4).Value = 67
$ws1.Cells.Item(333,5).Value = "LP1912"
$ws1.Cells.Item(334,1).Value = "16:45:34"
$ws1.Cells.Item(334,2).Value = "17:45"
$ws1.Cells.Item(334,3).Value = "15_ABASTO"
$ws1.Cells.Item(334,4).Value = 60
$ws1.Cells.Item(334,5).Value = "LP1912"
$ws1.Cells.Item(335,1).Value = "16:34:19"
$ws1.Cells.Item(335,2).Value = "17:50"
$ws1.Cells.Item(335,3).Value = "16_P MOR-167 Y 521"
$ws1.Cells.Item(335,4).Value = 76
$ws1.Cells.Item(335,5).Value = "LP1912"
$ws1.Cells.Item(336,1).Value = "15:59:02"
$ws1.Cells.Item(336,2).Value = "17:51"
$ws1.Cells.Item(336,3).Value = "16_P MOR-167 Y 521"
$ws1.Cells.Item(336,4).Value = 112
$ws1.Cells.Item(336,5).Value = "LP1912"
$ws1.Cells.Item(337,1).Value = "15:59:02"
$ws1.Cells.Item(337,2).Value = "17:52"
$ws1.Cells.Item(337,3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(337,4).Value = 113
$ws1.Cells.Item(337,5).Value = "LP1912"
$ws1.Cells.Item(338,1).Value = "16:20:15"
$ws1.Cells.Item(338,2).Value = "18:04"
$ws1.Cells.Item(338,3).Value = "17_ROMERO"
$ws1.Cells.Item(338,4).Value = 104
$ws1.Cells.Item(338,5).Value = "LP1912"
$ws1.Cells.Item(339,1).Value = "16:34:19"
$ws1.Cells.Item(339,2).Value = "18:21"
$ws1.Cells.Item(339,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(339,4).Value = 107
$ws1.Cells.Item(339,5).Value = "LP1912"
$ws1.Cells.Item(340,1).Value = "16:34:19"
$ws1.Cells.Item(340,2).Value = "18:28"
$ws1.Cells.Item(340,3).Value = "215C_EL PATO"
$ws1.Cells.Item(340,4).Value = 114
$ws1.Cells.Item(340,5).Value = "LP1912"
$ws1.Cells.Item(341,1).Value = "16:34:19"
$ws1.Cells.Item(341,2).Value = "18:32"
$ws1.Cells.Item(341,3).Value = "11X44_ETCHEVERRY"
$ws1.Cells.Item(341,4).Value = 118
$ws1.Cells.Item(341,5).Value = "LP1912"

# --- Sheet2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A2").Value = "Última actualización: 16:45:34"

# --- Sheet3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A2").Value = "Última actualización: 16:45:34"
$ws3.Range("A3").Value = "Total filas: 46"
$ws3.Rows.Item(50).Insert()
$ws3.Cells.Item(50,1).Value = "16:45:34"
$ws3.Cells.Item(50,2).Value = "18:03"
$ws3.Cells.Item(50,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(50,4).Value = 78
$ws3.Cells.Item(50,5).Value = "L6203"
